$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the Title/Company/City columns so the new sponsorship notes are easier to scan
$ws.Columns.Item(1).ColumnWidth = 49.333333333333336
$ws.Columns.Item(2).ColumnWidth = 45.5
$ws.Columns.Item(3).ColumnWidth = 32

# Remove stray empty City cells for rows 25 and 26 (filter leaves no city for these two postings)
$ws.Range("D25").ClearContents() | Out-Null
$ws.Range("D26").ClearContents() | Out-Null

# Append newly found sponsorship-friendly job postings (rows 41-48) with filter-result notes in column G
# Row 41
$ws.Range("A41").Value = "gXGih_3idtegO6LOAAAAAA=="
$ws.Range("B41").Value = "Entry Level QA (H1b Visa Sponsorship Available)"
$ws.Range("C41").Value = "Perfict Global, Inc."
$ws.Range("D41").Value = "New York"
$ws.Range("E41").Value = "New York"
$ws.Range("F41").Value = "https://www.optnation.com/entry-level-qa-h1b-visa-sponsorship-available-job-in-new-york-ny-view-jobid-33407?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G41").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 42
$ws.Range("A42").Value = "dh3Ad6iNO2ZMLcIBAAAAAA=="
$ws.Range("B42").Value = "Training&Placement in Business Analyst along with certification|VISA Sponsorship"
$ws.Range("C42").Value = "PrecisionTechnologies Corp"
$ws.Range("D42").Value = "South Brunswick Township"
$ws.Range("E42").Value = "New Jersey"
$ws.Range("F42").Value = "https://applicant.mightyrecruiter.com/jobs/apply/training-placement-in-business-analyst-along-with-certification-visa-sponsorship-at-precisiontechnologies-corp-in-south-brunswick-township-nj-e35475dc6941c5c2f182ca3f2182b448?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G42").Value = "Allows: explicit sponsorship signal ('VISA Sponsorship')"

# Row 43
$ws.Range("A43").Value = "pUeEVpkGzM8fAmDyAAAAAA=="
$ws.Range("B43").Value = "QA Analyst (H1b Visa Sponsorship Available)"
$ws.Range("C43").Value = "Perfict Global, Inc."
$ws.Range("D43").Value = "Boston"
$ws.Range("E43").Value = "Massachusetts"
$ws.Range("F43").Value = "https://www.optnation.com/qa-analyst-h1b-visa-sponsorship-available-job-in-boston-ma-view-jobid-33582?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G43").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 44
$ws.Range("A44").Value = "eZczUsY6W1vQHTaGAAAAAA=="
$ws.Range("B44").Value = "Entry Level QA (H1b Visa Sponsorship Available)"
$ws.Range("C44").Value = "Perfict Global, Inc."
$ws.Range("D44").Value = "Pittsburgh"
$ws.Range("E44").Value = "Pennsylvania"
$ws.Range("F44").Value = "https://www.optnation.com/entry-level-qa-h1b-visa-sponsorship-available-job-in-pittsburgh-pa-view-jobid-33580?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G44").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 45
$ws.Range("A45").Value = "xoM74PXhUYkigicIAAAAAA=="
$ws.Range("B45").Value = "Entry Level QA (H1b Visa Sponsorship Available)"
$ws.Range("C45").Value = "Perfict Global, Inc."
$ws.Range("D45").Value = "Louisville"
$ws.Range("E45").Value = "Kentucky"
$ws.Range("F45").Value = "https://www.optnation.com/entry-level-qa-h1b-visa-sponsorship-available-job-in-louisville-ky-view-jobid-33765?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G45").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 46
$ws.Range("A46").Value = "yv24shpMgCVxcssFAAAAAA=="
$ws.Range("B46").Value = "Junior Software Developer – USA Visa Sponsorship in USA – (job id: 1681772381)"
$ws.Range("C46").Value = "vmysmartpros"
$ws.Range("F46").Value = "https://www.mysmartpros.com/tuition/job/junior-software-developer-usa-visa-sponsorship-in-usa-job-id-1681772381/?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G46").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 47
$ws.Range("A47").Value = "4vszmTyq5cJb08j9AAAAAA=="
$ws.Range("B47").Value = "Entry Level Software Tester (H1b Visa Sponsorship available)"
$ws.Range("C47").Value = "athomejobs5.10001mb"
$ws.Range("F47").Value = "https://athomejobs5.10001mb.com/job/entry-level-software-tester-h1b-visa-sponsorship-available?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G47").Value = "Allows: explicit sponsorship signal ('Visa Sponsorship')"

# Row 48
$ws.Range("A48").Value = "KYU9wwLCb_wuuemFAAAAAA=="
$ws.Range("B48").Value = "Data Engineer at BeaconFire Inc. Trenton, NJ"
$ws.Range("C48").Value = "BeaconFire Inc."
$ws.Range("D48").Value = "Trenton"
$ws.Range("E48").Value = "New Jersey"
$ws.Range("F48").Value = "https://badidearadio.com/job-library/job/data-engineer-at-beaconfire-inc-trenton-nj-a05vcmNiRC9JbEl1VkJJelhXYWN1aHpUeVE9PQ==?utm_campaign=google_jobs_apply&utm_source=google_jobs_apply&utm_medium=organic"
$ws.Range("G48").Value = "Allows: explicit sponsorship signal ('visa sponsorship')"

# Leave the selection where the review of the newly-added rows left off
$ws.Range("F24").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 12
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll-position automation isn't exposed everywhere; selection above already covers the cursor.
}
